$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.968.14"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.556.67"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "207.37"
$ws.Range("E6").Value = "  +0.61%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.16"
$ws.Range("E8").Value = "  +4.18%  "
$ws.Range("E9").Value = "  +0.29%  "
$ws.Range("E10").Value = "  +1.15%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0857"
$ws.Range("E11").Value = "  +0.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.778.94"
$ws.Range("E12").Value = "  +0.78%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.556.74"
$ws.Range("E13").Value = "  +0.73%  "
$ws.Range("E14").Value = "  +1.53%  "
$ws.Range("E15").Value = "  +1.88%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.970.79"
$ws.Range("E16").Value = "  +0.37%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.82"
$ws.Range("E17").Value = "  +0.40%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "218.49"
$ws.Range("E18").Value = "  +2.28%  "
$ws.Range("E19").Value = "  +2.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.32"
$ws.Range("E20").Value = "  +2.01%  "
$ws.Range("E21").Value = "  -0.16%  "
$ws.Range("E22").Value = "  +1.22%  "
$ws.Range("E23").Value = "  +0.90%  "
$ws.Range("E24").Value = "  +0.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.19"
$ws.Range("E25").Value = "  +0.99%  "
$ws.Range("E26").Value = "  +0.61%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.98"
$ws.Range("E27").Value = "  +1.20%  "
$ws.Range("E28").Value = "  +1.20%  "
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("E30").Value = "  +2.49%  "
$ws.Range("E31").Value = "  -0.17%  "
$ws.Range("E32").Value = "  +0.84%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.425.22"
$ws.Range("E33").Value = "  +5.11%  "
$ws.Range("E34").Value = "  +5.16%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.59"
$ws.Range("E35").Value = "  +3.99%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.982"
$ws.Range("E36").Value = "  +1.95%  "
$ws.Range("E37").Value = "  +0.23%  "
$ws.Range("E38").Value = "  +0.86%  "
$ws.Range("E39").Value = "  +0.34%  "
$ws.Range("E40").Value = "  +1.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.76"
$ws.Range("E41").Value = "  +3.08%  "
$ws.Range("E42").Value = "  -0.12%  "
$ws.Range("E43").Value = "  +5.17%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.986"
$ws.Range("E44").Value = "  -0.37%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.52"
$ws.Range("E45").Value = "  +1.78%  "
$ws.Range("E46").Value = "  +2.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.692.44"
$ws.Range("E47").Value = "  +0.77%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.00"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0521"
$ws.Range("E49").Value = "  +2.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₇0998"
$ws.Range("E50").Value = "  +3.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0957"
$ws.Range("E51").Value = "  +1.07%  "
